# Fill in the risk table ("Riscos") of the Scrum risk-management workbook
# with the set of risks identified in the latest review.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Riscos")

# Probability (E) / Impact (F) numbers, by row
$ws.Cells.Item(3, 5).Value = 3
$ws.Cells.Item(3, 6).Value = 2

$ws.Cells.Item(4, 5).Value = 1
$ws.Cells.Item(4, 6).Value = 2

$ws.Cells.Item(5, 5).Value = 1
$ws.Cells.Item(5, 6).Value = 2

$ws.Cells.Item(6, 5).Value = 1
$ws.Cells.Item(6, 6).Value = 3

$ws.Cells.Item(7, 5).Value = 3
$ws.Cells.Item(7, 6).Value = 3

$ws.Cells.Item(8, 5).Value = 1
$ws.Cells.Item(8, 6).Value = 1

$ws.Cells.Item(9, 5).Value = 2
$ws.Cells.Item(9, 6).Value = 3

$ws.Cells.Item(10, 5).Value = 3
$ws.Cells.Item(10, 6).Value = 1

$ws.Cells.Item(11, 5).Value = 2
$ws.Cells.Item(11, 6).Value = 3

$ws.Cells.Item(12, 5).Value = 3
$ws.Cells.Item(12, 6).Value = 3

# Risk description (D) / impact description (G) text, entered in the same
# order the original author typed them in.
$ws.Cells.Item(3, 4).Value = "Atraso da Entrega da Sprint"
$ws.Cells.Item(4, 4).Value = "Estouro do Orçamento"
$ws.Cells.Item(4, 7).Value = "O estouro no orçamento acaba culminando em um replanejamento do projeto e das tarefas seguintes."
$ws.Cells.Item(3, 7).Value = "O atraso da entrega pode comprometer todo o cronograma previamente estipulado."

$ws.Cells.Item(5, 4).Value = "Desistência de Algum Membro do Projeto"
$ws.Cells.Item(5, 7).Value = "Caso algum membro desista os outros devem dividir o que o desistente iria fazer, e replanejar as horas dedicada."

$ws.Cells.Item(6, 4).Value = "Resultado Desaprovado pelo PO"
$ws.Cells.Item(6, 7).Value = "Caso o PO desaprove o resultado da sprint, todo o projeto deve ser replanejado, incluindo datas."

$ws.Cells.Item(7, 4).Value = "Equipe Inexperiente"
$ws.Cells.Item(8, 4).Value = "Falta de treinamento com Ferramentas"
$ws.Cells.Item(7, 7).Value = "Por ser uma equipe nova, a inexperiência irá alocar mais tempo para aprendizado."
$ws.Cells.Item(8, 7).Value = "Algumas ferramentas utilizadas necessitam de aprendizado pela equipe."

$ws.Cells.Item(9, 4).Value = "Documentação Atrasada"
$ws.Cells.Item(9, 7).Value = "Se a documentação do projeto atrasar, não será possível prosseguir com outras áreas."

$ws.Cells.Item(10, 4).Value = "Impossibilidade de Reunião"
$ws.Cells.Item(10, 7).Value = "Algum membro pode não estar disponível para qualquer das reuniões"

$ws.Cells.Item(11, 4).Value = "Quantidade Alta de Defeitos"
$ws.Cells.Item(11, 7).Value = "Os defeitos no software a ser entregue pode ser grande, comprometendo a entrega."

$ws.Cells.Item(12, 4).Value = "Não Cumprimento de Tarefas"
$ws.Cells.Item(12, 7).Value = "Algum membro pode não cumprir o que for alocado para o mesmo realizar."

# Rows with longer descriptions wrap to more text lines; reflect the taller
# row heights that Excel computes once the wrapped text is present.
$ws.Rows.Item(3).RowHeight = 60
$ws.Rows.Item(4).RowHeight = 60
$ws.Rows.Item(5).RowHeight = 60
$ws.Rows.Item(6).RowHeight = 60
$ws.Rows.Item(7).RowHeight = 45
$ws.Rows.Item(8).RowHeight = 45
$ws.Rows.Item(9).RowHeight = 45
$ws.Rows.Item(10).RowHeight = 45
$ws.Rows.Item(11).RowHeight = 45
$ws.Rows.Item(12).RowHeight = 45

# Update the active selection to reflect where the editor left off
$ws.Range("M9").Select() | Out-Null

$wb.Save()
